# "fixed some data and add some cuts"
# Update measured intensity values for a couple of frames; the dependent
# formulas (G, H, I, J columns) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Fotograma 1): new I_max / I_min1 / I_min2 readings
$ws.Range("D2").Value = 185.1
$ws.Range("E2").Value = 142.1
$ws.Range("F2").Value = 146.7

# Row 8 (Fotograma 7): new I_max / I_min1 / I_min2 readings
$ws.Range("D8").Value = 137.5
$ws.Range("E8").Value = 112.5
$ws.Range("F8").Value = 118.1

# Move the active selection to I8, matching the author's last cursor spot
$ws.Range("I8").Select()
